# refactor: merapikan kode impor file
# - cell A4 (id_ruangan for the "AC Panasonic" facility row) changes from 3 to 0
# - active selection moves from A6 to A4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 0

$ws.Range("A4").Select()
